$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.586.39"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "2.110.72"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  +0.91%  "
$ws.Range("D5").Value = "'336.40"
$ws.Range("E6").Value = "  +0.88%  "
$ws.Range("E7").Value = "  +0.99%  "
$ws.Range("D8").Value = "'0.4549"
$ws.Range("E8").Value = "  +4.85%  "
$ws.Range("D9").Value = "'55.47"
$ws.Range("E9").Value = "  +4.85%  "
$ws.Range("D10").Value = "'0.08995"
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("E11").Value = "  +1.34%  "
$ws.Range("D12").Value = "'24.69"
$ws.Range("E12").Value = "  +1.48%  "
$ws.Range("D13").Value = "2.111.78"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").Value = "'6.861"
$ws.Range("E14").Value = "  +2.81%  "
$ws.Range("D15").Value = "'8.128"
$ws.Range("E15").Value = "  +6.05%  "
$ws.Range("E16").Value = "  +5.21%  "
$ws.Range("D17").Value = "'97.33"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("E18").Value = "  +0.93%  "
$ws.Range("D19").Value = "'0.06689"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("D20").Value = "'19.35"
$ws.Range("E20").Value = "  +0.79%  "
$ws.Range("E21").Value = "  +0.93%  "
$ws.Range("D22").Value = "'6.260"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "30.650.69"
$ws.Range("E23").Value = "  +0.52%  "
$ws.Range("D24").Value = "'12.79"
$ws.Range("E24").Value = "  +5.12%  "
$ws.Range("D25").Value = "'2.363"
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("D26").Value = "2.361.10"
$ws.Range("E26").Value = "  +1.25%  "
$ws.Range("D27").Value = "'22.30"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("D28").Value = "'163.82"
$ws.Range("E28").Value = "  +0.87%  "
$ws.Range("D29").Value = "'2.519"
$ws.Range("E29").Value = "  -2.38%  "
$ws.Range("D30").Value = "'133.61"
$ws.Range("E30").Value = "  +1.67%  "
$ws.Range("D31").Value = "'1.222"
$ws.Range("E31").Value = "  +2.79%  "
$ws.Range("D32").Value = "'0.1071"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").Value = "'6.361"
$ws.Range("E33").Value = "  +3.96%  "
$ws.Range("D34").Value = "'1.628"
$ws.Range("D35").Value = "'3.967"
$ws.Range("E35").Value = "  +2.12%  "
$ws.Range("D36").Value = "'10.50"
$ws.Range("E36").Value = "  +4.24%  "
$ws.Range("D37").Value = "'5.870"
$ws.Range("E37").Value = "  +7.82%  "
$ws.Range("D38").Value = "'0.02612"
$ws.Range("E38").Value = "  +1.63%  "
$ws.Range("D39").Value = "'0.06824"
$ws.Range("E39").Value = "  +0.46%  "
$ws.Range("D40").Value = "'0.2317"
$ws.Range("E40").Value = "  +2.63%  "
$ws.Range("E41").Value = "  -0.98%  "
$ws.Range("D42").Value = "'0.6854"
$ws.Range("E42").Value = "  -0.53%  "
$ws.Range("D43").Value = "'1.255"
$ws.Range("E43").Value = "  -0.47%  "
$ws.Range("D44").Value = "'0.6446"
$ws.Range("E44").Value = "  +1.42%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'14.12"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "'2.313"
$ws.Range("E46").Value = "  +5.28%  "
$ws.Range("E47").Value = "  +1.64%  "
$ws.Range("D48").Value = "'1.251"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").Value = "'0.00000000346"
$ws.Range("E49").Value = "  +17.77%  "
$ws.Range("D50").Value = "'1.211"
$ws.Range("E50").Value = "  -2.03%  "
$ws.Range("D51").Value = "'83.15"
$ws.Range("E51").Value = "  +1.65%  "
